$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.231.31"
$ws.Range("E2").Value = "  -0.77%  "

$ws.Range("D3").Value = "'1.841.55"
$ws.Range("E3").Value = "  -0.54%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'240.59"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("D6").Value = "'0.6272"
$ws.Range("E6").Value = "  -0.28%  "

$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").Value = "'0.07442"
$ws.Range("E8").Value = "  -2.81%  "

$ws.Range("D9").Value = "'0.2885"
$ws.Range("E9").Value = "  -1.16%  "

$ws.Range("D10").Value = "'24.23"
$ws.Range("E10").Value = "  -2.39%  "

$ws.Range("D11").Value = "'0.07734"
$ws.Range("E11").Value = "  -0.13%  "

$ws.Range("D12").Value = "'1.840.70"
$ws.Range("E12").Value = "  -2.52%  "

$ws.Range("D13").Value = "'4.983"
$ws.Range("E13").Value = "  -1.08%  "

$ws.Range("D14").Value = "'0.6755"
$ws.Range("E14").Value = "  -0.88%  "

$ws.Range("D15").Value = "'0.00001009"
$ws.Range("E15").Value = "  -4.12%  "

$ws.Range("E16").Value = "  -1.82%  "

$ws.Range("D17").Value = "'6.110"
$ws.Range("E17").Value = "  -1.38%  "

$ws.Range("D18").Value = "'29.277.68"

$ws.Range("D19").Value = "'227.23"
$ws.Range("E19").Value = "  -0.87%  "

$ws.Range("D20").Value = "'12.25"
$ws.Range("E20").Value = "  -0.76%  "

$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  +0.16%  "

$ws.Range("D22").Value = "'7.340"
$ws.Range("E22").Value = "  -1.86%  "

$ws.Range("D23").Value = "'1.003"
$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").Value = "'158.77"
$ws.Range("E24").Value = "  +0.77%  "

$ws.Range("D25").Value = "'0.1370"
$ws.Range("E25").Value = "  -1.01%  "

$ws.Range("D26").Value = "'8.371"
$ws.Range("E26").Value = "  -0.69%  "

$ws.Range("D27").Value = "'17.53"
$ws.Range("E27").Value = "  -1.34%  "

$ws.Range("E28").Value = "  +1.57%  "

$ws.Range("D29").Value = "'0.06103"
$ws.Range("E29").Value = "  +8.70%  "

$ws.Range("D30").Value = "'1.470"
$ws.Range("E30").Value = "  +0.32%  "

$ws.Range("D31").Value = "'4.071"
$ws.Range("E31").Value = "  -1.60%  "

$ws.Range("D32").Value = "'4.031"
$ws.Range("E32").Value = "  -0.88%  "

$ws.Range("D33").Value = "'1.816"
$ws.Range("E33").Value = "  -1.58%  "

$ws.Range("D34").Value = "'1.138"
$ws.Range("E34").Value = "  -2.56%  "

$ws.Range("D35").Value = "'0.6955"
$ws.Range("E35").Value = "  -0.74%  "

$ws.Range("D36").Value = "'2.594"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").Value = "'2.822"
$ws.Range("E37").Value = "  +2.98%  "

$ws.Range("D38").Value = "'1.241.57"
$ws.Range("E38").Value = "  +1.15%  "

$ws.Range("D39").Value = "'0.01811"
$ws.Range("E39").Value = "  +0.38%  "

$ws.Range("D40").Value = "'6.491"
$ws.Range("E40").Value = "  +0.33%  "

$ws.Range("D41").Value = "'0.9073"
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  -0.09%  "

$ws.Range("D43").Value = "'1.999.61"
$ws.Range("E43").Value = "  -5.07%  "

$ws.Range("D44").Value = "'101.49"
$ws.Range("E44").Value = "  -0.67%  "

$ws.Range("D45").Value = "'66.00"
$ws.Range("E45").Value = "  -0.26%  "

$ws.Range("D46").Value = "'7.026"
$ws.Range("E46").Value = "  -2.80%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.1161"
$ws.Range("E47").Value = "  +0.42%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.981"
$ws.Range("E48").Value = "  -0.85%  "

$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").Value = "'0.3931"
$ws.Range("E49").Value = "  -2.38%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.00000000114"
$ws.Range("E50").Value = "  -3.22%  "

$ws.Range("D51").Value = "'1.652"
$ws.Range("E51").Value = "  -1.63%  "
